$wb = $excel.ActiveWorkbook

# --- Sheet "NP 10": add a new pDNA measurement block in row 27 ---
$npTen = $wb.Worksheets.Item("NP 10")
$npTen.Range("D27").Value = 252.77
$npTen.Range("E27").Value = 266.12
$npTen.Range("F27").Value = 290.02
$npTen.Range("G27").Formula = "=AVERAGE(D27:F27)"
$npTen.Range("H27").Formula = "=STDEV(D27:F27)"
$npTen.Range("D27:H27").Style = "Normal"

# --- Sheet "Sheet1": insert a new summary row for the pDNA group ---
$summary = $wb.Worksheets.Item("Sheet1")
$summary.Rows.Item(2).Insert()
$summary.Range("A2").Value = "pDNA"
$summary.Range("B2").Value = 0
$summary.Range("C2").Formula = "=269.636666666667/2"
$summary.Range("D2").Formula = "=18.872356327002/2"
$summary.Range("A2:D2").Style = "Normal"

# --- Restore / update the view selections shown in the diff ---
[void]$npTen.Range("G27:H27").Select()
[void]$summary.Range("F23").Select()
